$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.894424333333333
$ws.Range("H2").Value = 5.683273
$ws.Range("I2").Value = 0.6002819911800915
$ws.Range("J2").Value = 0.6002819911800916
$ws.Range("M2").Value = 0.04554333333333333
$ws.Range("N2").Value = 0.13663
$ws.Range("Q2").Value = 0.08627839888777777
$ws.Range("R2").Value = 0.77650558999
$ws.Range("S2").Value = 0.6002819911800915
$ws.Range("T2").Value = 0.6002819911800916

# Row 3
$ws.Range("I3").Value = 0.3198928944728968
$ws.Range("J3").Value = 0.3198928944728969
$ws.Range("M3").Value = 0.04554333333333333
$ws.Range("N3").Value = 0.13663
$ws.Range("Q3").Value = 0.04597813553666667
$ws.Range("R3").Value = 0.4138032198300001
$ws.Range("S3").Value = 0.3198928944728968
$ws.Range("T3").Value = 0.3198928944728969

# Row 4
$ws.Range("G4").Value = 0.1798433333333334
$ws.Range("H4").Value = 0.5395300000000001
$ws.Range("I4").Value = 0.05698655382231239
$ws.Range("J4").Value = 0.05698655382231241
$ws.Range("M4").Value = 0.04554333333333333
$ws.Range("N4").Value = 0.13663
$ws.Range("Q4").Value = 0.008190664877777778
$ws.Range("R4").Value = 0.0737159839
$ws.Range("S4").Value = 0.05698655382231239
$ws.Range("T4").Value = 0.05698655382231241

# Row 5
$ws.Range("G5").Value = 0.07207599999999999
$ws.Range("H5").Value = 0.216228
$ws.Range("I5").Value = 0.0228385605246992
$ws.Range("J5").Value = 0.02283856052469921
$ws.Range("M5").Value = 0.04554333333333333
$ws.Range("N5").Value = 0.13663
$ws.Range("Q5").Value = 0.003282581293333333
$ws.Range("R5").Value = 0.02954323164
$ws.Range("S5").Value = 0.0228385605246992
$ws.Range("T5").Value = 0.02283856052469921
